# Fix formatting of floating point "Importe" values that were scraped with
# European-style thousands/decimal separators (e.g. "1.404,92") but need to
# read as plain decimal numbers formatted with a period decimal point and no
# thousands separator (e.g. "1404.92"). Also fix a stray comma that should be
# a period in a provider name in column E.
#
# The target cells are stored as literal text (shared strings) in the
# workbook, not as real numbers, so a straight `$cell.Value = "..."` would
# get auto-coerced by Excel into a genuine number (dropping formatting, e.g.
# "510.00" -> 510). To keep the result as literal text without touching the
# number format of the destination cells themselves, we stage each new
# string in a scratch cell that IS formatted as Text, copy it, and use
# PasteSpecial (values only) into the destination - this carries over the
# literal text without reformatting the destination cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the stray comma in the provider name (plain text, not numeric-looking,
# so a normal assignment is safe and keeps it as text).
$ws.Range("E59").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"

# Map of cell -> corrected "Importe" text.
$importeFixes = @{
    "H2" = "510.00"
    "H3" = "46320.00"
    "H4" = "1404.92"
    "H5" = "100.80"
    "H6" = "5714.37"
    "H7" = "62973.20"
    "H8" = "12600.46"
    "H9" = "146.65"
    "H10" = "571.60"
    "H11" = "7214.19"
    "H12" = "300.00"
    "H13" = "1790.00"
    "H14" = "4615.20"
    "H15" = "22.00"
    "H16" = "29998.78"
    "H17" = "861.60"
    "H18" = "15.00"
    "H19" = "3057.00"
    "H20" = "1400.00"
    "H21" = "2300.00"
    "H22" = "67.25"
    "H23" = "2700.00"
    "H24" = "4677.40"
    "H25" = "268.00"
    "H26" = "129.60"
    "H27" = "30.00"
    "H28" = "1.22"
    "H29" = "67500.00"
    "H30" = "67500.00"
    "H31" = "67500.00"
    "H32" = "67500.00"
    "H33" = "67500.00"
    "H34" = "103769.60"
    "H35" = "9.00"
    "H36" = "1817.25"
    "H37" = "325.00"
    "H38" = "149.00"
    "H39" = "225.16"
    "H40" = "192.00"
    "H41" = "10240.00"
    "H42" = "4235.00"
    "H43" = "2230.00"
    "H44" = "1800.00"
    "H45" = "239309.00"
    "H46" = "238142.20"
    "H47" = "500.00"
    "H48" = "700.00"
    "H49" = "580.00"
    "H50" = "250.00"
    "H51" = "1657.50"
    "H52" = "200.00"
    "H53" = "3400.00"
    "H54" = "750.00"
    "H55" = "120.00"
    "H56" = "1000.00"
    "H57" = "4400.00"
    "H58" = "99.00"
    "H59" = "170.00"
    "H60" = "225.00"
    "H61" = "1930.00"
    "H62" = "1182.40"
    "H63" = "29727.68"
    "H64" = "2690.10"
    "H65" = "546.00"
    "H66" = "65000.00"
    "H67" = "44000.00"
    "H68" = "878.35"
    "H69" = "2000.00"
    "H70" = "3978.00"
    "H71" = "3626.00"
    "H72" = "2155000.00"
    "H73" = "800.00"
}

# Scratch cell well outside the used range, formatted as Text so values
# typed/pasted into it are never reinterpreted as numbers.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

foreach ($cellRef in $importeFixes.Keys) {
    $scratch.Value = $importeFixes[$cellRef]
    $scratch.Copy()
    $dest = $ws.Range($cellRef)
    $dest.PasteSpecial(-4163)
}

# Remove the scratch cell's content/format so it doesn't linger in the
# workbook's used range.
$scratch.Clear()

$excel.CutCopyMode = 0
